# TC10_CDS_Filter_PHSAccession-phs002504.xlsx correction
# Remove the "RNA-Seq" experimental_strategies filter from the FilesTab
# query (row 4, column B) so the list is empty again, matching the other
# tabs' queries. Also restores the selection state left behind after the
# edit (Excel moves the active cell while the correction is made).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ("FilesTab") query correction -----------------------------
$filesTabCell = $ws.Cells.Item(4, 2)
$query = $filesTabCell.Value()
$fixedQuery = $query.Replace('experimental_strategies: ["RNA-Seq"]', 'experimental_strategies: []')
$filesTabCell.Value = $fixedQuery

# Keep the (very tall, word-wrapped) rows at their previous rendered
# height now that the text is a little shorter.
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# --- Restore view/selection -------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 2
$ws.Range("C4").Select() | Out-Null
